# Change "Hearing Type" column header/values to "Type" on the "Hearing List" sheet,
# then leave "Hearing List" as the active sheet with D14 selected.

$wb = $excel.ActiveWorkbook

$wsHearing = $wb.Worksheets.Item("Hearing List")

# Rename the header and the values in column D ("Hearing Type" -> "Type")
$wsHearing.Range("D1").Value = "Type"
$wsHearing.Range("D2").Value = "Type A"
$wsHearing.Range("D3").Value = "Type B"
$wsHearing.Range("D4").Value = "Type C"

# Make "Hearing List" the active sheet, and select D14 on it.
$wsHearing.Activate()
$wsHearing.Range("D14").Select()
